# Sync attendance_reports: reorder "Recorded By" (column G) values so that
# the leading "System, " token is moved to the end of the comma-separated
# list instead of the front, e.g.
#   "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"
#   "System, system, backup@backdoor.com" -> "system, backup@backdoor.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$prefix = "System, "

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $value = $cell.Value2

    if ($null -ne $value -and $value -is [string] -and $value.StartsWith($prefix)) {
        $rest = $value.Substring($prefix.Length)
        $newValue = $rest + ", System"
        $cell.Value2 = $newValue
    }
}
